$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new entry dated 5 February 2021
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B4").Value = "Implemented Feed Forward Neural Network model and Convolutional Neural Network. The Feed Forward Neural Network achieved final accuracy of 88.69% on Fruits-360 dataset. Whereas, Convolutional Neural Network achieved accuracy of 95% on Fruits-360 dataset."

$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C4").Value = "Apply different algorithms to the Fruits-360 dataset and also combine Fruits-360 dataset with another dataset (waste), to see the results."

$excel.CutCopyMode = 0

# Date serial number for 5 February 2021 (avoids the engine auto-assigning
# a datetime display format when given a DateTime object directly)
$ws.Range("A4").Value = 44232
$ws.Range("A4").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A4").VerticalAlignment = -4108    # xlCenter
$ws.Range("A4").WrapText = $true
$ws.Range("A4").NumberFormat = "mm-dd-yy"

$ws.Rows.Item(4).RowHeight = 100.8

# Update selection to reflect the newly active cell
$ws.Range("C4").Select()
